# Update cryptocurrency price/volume data per upstream refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.412.52"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.13%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.514.80"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -2.32%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.04%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'606.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.51%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'143.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -3.94%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.514.68"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -2.29%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -0.22%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.511"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +4.63%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = "'Toncoin"
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = "'7.71"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -4.64%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = "'Dogecoin"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = "'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = "'0.130"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -4.48%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.407"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -2.08%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'4.127.51"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -2.14%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  -7.26%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'28.60"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -4.09%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.522.87"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -2.57%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  +0.66%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'66.321.96"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.46%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'10.80"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -6.22%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'6.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -3.92%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'14.58"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -3.56%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'421.90"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -1.39%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.588"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -5.02%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'76.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -2.30%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'3.663.42"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -2.32%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D27").Value = "'0.0000113"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -7.80%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'7.88"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -5.41%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -2.53%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'8.88"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -5.77%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'1.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +0.06%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'3.528.49"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -2.01%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -1.72%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'24.11"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -5.39%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  +0.03%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'1.33"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -9.70%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'7.52"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -4.41%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -4.91%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'173.75"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -2.03%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'5.18"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -8.77%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.0815"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -4.87%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'4.96"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -5.30%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.855"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -4.96%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'45.33"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -1.39%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'1.75"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -8.28%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'1.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.04%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -8.78%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'7.04"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -2.29%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'1.10"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -6.92%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'22.65"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -5.92%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.901"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -5.55%  "
$ws.Range("E51").Style = "Normal"
